# Updated main GSC export data.
#
# The "Chart" sheet's first data row (2025-11-02, which only ever held
# placeholder/blank indexing counts) is dropped entirely; every following
# row shifts up one position so each date keeps its own indexing figures.
# After the shift, the sheet's final row (now 2026-01-29) gets a proper
# numeric 0 for Impressions instead of the leftover blank/text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete 2025-11-02 placeholder row; Excel shifts rows 3:90 up to 2:89.
$ws.Rows.Item(2).Delete()

# The row that is now last (2026-01-29) should carry a numeric Impressions value.
$lastRow = $ws.Cells(1, 1).End(-4121).Row
$ws.Cells.Item($lastRow, 4).Value = 0
